# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback packages have been generated: the Status moves
# from "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns are filled in for each language table, and column widths are
# widened to fit the new content.

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/3d8fb30f10e9d9cb967a4bd3808eb1e7baa7832c/e2e/"

$mdFile1 = "b243efce-f9ab-49c9-b4c8-f186f6f05ffa.md"
$mdFile2 = "fb3af559-1345-4472-89ad-fa43174ee67a.md"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status cells
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("K2").Value = "2016-08-12 22:35:45"
$wsZh.Range("K3").Value = "2016-08-12 22:35:45"

$wsZh.Range("J2").Value = "b243efce-f9ab-49c9-b4c8-f186f6f05ffa.4d0c3b5a2f326cc091a133b0fba0517d376bf8f3.zh-cn.xlf"
$wsZh.Range("J3").Value = "fb3af559-1345-4472-89ad-fa43174ee67a.a4b5207318d5c4283ba47cdab6ea3db81e1c84b4.zh-cn.xlf"

# Rebuild the hyperlinks in order A2, I2, A3, I3 so relationship ids come
# out sequential (rId2..rId5), matching the new "Latest Target File" links.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $ghBase + $mdFile1, "", "", $mdFile1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $ghBase + $mdFile1, "", "", $mdFile1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $ghBase + $mdFile2, "", "", $mdFile2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $ghBase + $mdFile2, "", "", $mdFile2) | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("K2").Value = "2016-08-12 22:35:55"
$wsDe.Range("K3").Value = "2016-08-12 22:35:55"

$wsDe.Range("J2").Value = "b243efce-f9ab-49c9-b4c8-f186f6f05ffa.4d0c3b5a2f326cc091a133b0fba0517d376bf8f3.de-de.xlf"
$wsDe.Range("J3").Value = "fb3af559-1345-4472-89ad-fa43174ee67a.a4b5207318d5c4283ba47cdab6ea3db81e1c84b4.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $ghBase + $mdFile1, "", "", $mdFile1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $ghBase + $mdFile1, "", "", $mdFile1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $ghBase + $mdFile2, "", "", $mdFile2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $ghBase + $mdFile2, "", "", $mdFile2) | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated for zh-cn and de-de."
